$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (bold, border, alignment => style index 1) from an
# existing header cell onto the three new header cells, then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every data row (2-41).
$ws.Range("AD2:AD41").Value = 98
$ws.Range("AE2:AE41").Value = 64
$ws.Range("AF2:AF41").Value = 0
